$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "322×8=2576" "841×7=5887"
Replace-Text "495×4=1980" "362×4=1448"
Replace-Text "358×6=2148" "711×2=1422"
Replace-Text "381×8=3048" "191×9=1719"
Replace-Text "789×3=2367" "273×5=1365"
Replace-Text "738×4=2952" "411×4=1644"
Replace-Text "156×2=312" "245×5=1225"
Replace-Text "426×5=2130" "448×9=4032"
Replace-Text "587×2=1174" "437×5=2185"
Replace-Text "140×2=280" "490×3=1470"
Replace-Text "516×2=1032" "683×2=1366"
Replace-Text "485×3=1455" "522×8=4176"
Replace-Text "848×4=3392" "149×2=298"
Replace-Text "259×9=2331" "169×8=1352"
Replace-Text "630×9=5670" "888×6=5328"
Replace-Text "529×2=1058" "261×3=783"
Replace-Text "451×6=2706" "836×6=5016"
Replace-Text "540×4=2160" "480×2=960"
Replace-Text "681×3=2043" "750×3=2250"
Replace-Text "774×8=6192" "344×9=3096"
Replace-Text "487×7=3409" "163×2=326"
Replace-Text "390×8=3120" "935×7=6545"
Replace-Text "972×4=3888" "874×9=7866"
Replace-Text "949×5=4745" "945×6=5670"
Replace-Text "496×7=3472" "885×2=1770"

Write-Host "Done applying replacements"
